$wb = $excel.ActiveWorkbook

# Update the "Status" cell text everywhere it is "Ready for handoff" -> "Handback transform failed"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Add the new "Error Detail" column (K) values describing the handback/handoff file name mismatch
$wsZhCn.Range("K3").Value = "Handback file name: i1etjgw3.qrr is different with handoff file name: 1610f715-4a0a-43da-8c92-b72a430f2de6.e5ea5b5364c7116a96feb5797fa9e663a41fdaa8.zh-cn."
$wsDeDe.Range("K3").Value = "Handback file name: i1etjgw3.qrr is different with handoff file name: 1610f715-4a0a-43da-8c92-b72a430f2de6.e5ea5b5364c7116a96feb5797fa9e663a41fdaa8.de-de."
